$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
questions = [
    {
        "title": "You want to modify your browser's behavior using the options parameter of the web driver.  Which of the following options should you import to do this with the Chrome browser?",
        "ques_type": 2,
        "options": [
            "from selenium.webdriver.chrome import ChromeOptions",
            "from selenium.webdriver.chrome.options import ChromiumOptions",
            "from selenium.webdriver.chrome.options import Options",
            "from selenium.webdriver.common.options import ArgOptions"
        ],
        "score": "from selenium.webdriver.chrome.options import Options"
    },
    {
        "title": "True or false: You can import a web driver object and use it to navigate to Google using the code shown below. from selenium import webdriver\nwebdriver.Chrome().navigate('https://google.com')",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    },
    {
        "title": "You\u2019re testing a new site, and you\u2019ve noticed that when a button on the page is covered by a pop-up, you get an exception when you try to click it.  Which of the following exceptions do you need to catch in this case?",
        "ques_type": 2,
        "options": [
            "ElementNotInteractableException",
            "NoSuchElementException",
            "ElementNotVisibleException",
            "ElementNotSelectableException"
        ],
        "score": "ElementNotInteractableException"
    },
    {
        "title": "You are working with a button on a webpage with the ID submit_form.  Which of the following methods is most appropriate to find and click this button?",
        "ques_type": 2,
        "options": [
            "driver.find_elements_by_css_selector('#submit_form').click()",
            "driver.find_element(By.CSS_SELECTOR, \"#submit_form\").click()",
            "driver.find_element(By.ID, 'submit_form').click()",
            "driver.find_elements_by_css_selector('#submit_form')[0].click()"
        ],
        "score": "driver.find_element(By.CSS_SELECTOR, \"#submit_form\").click()"
    }
]
'@

# Remove row 2 (the old un-pretty-printed duplicate) entirely, shifting nothing below it up.
$ws.Range("A2").EntireRow.Delete() | Out-Null

# A1 previously held a bare 0 with a bold/centered/bordered style; the new content
# is the shared-string text with plain/default formatting, so clear any formatting first.
$ws.Range("A1").ClearFormats() | Out-Null
$ws.Range("A1").Value = $newText

# Setting a multi-segment string can auto-expand the row height; restore natural auto height.
$ws.Rows(1).AutoFit() | Out-Null
